# The author fixed a grammar slip in section 1.2 "Objectives": the
# original sentence read "...is to create a intuitive mobile App..."
# (missing the "n" that turns "a" into "an" before the vowel sound in
# "intuitive"). Use Find/Replace over the document's logical text (this
# spans several runs / a couple of now-removed proofing-error markers in
# the underlying XML, which Find.Execute abstracts away).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$wdFindContinue = 1
$wdReplaceOne   = 1

$ok = $find.Execute(
    "create a intuitive mobile App",   # FindText
    $true,                              # MatchCase
    $false,                             # MatchWholeWord
    $false,                             # MatchWildcards
    $false,                             # MatchSoundsLike
    $false,                             # MatchAllWordForms
    $true,                              # Forward
    $wdFindContinue,                    # Wrap
    $false,                             # Format
    "create an intuitive mobile App",   # ReplaceWith
    $wdReplaceOne                       # Replace
)

if (-not $ok) {
    throw "Could not find the sentence to correct ('a intuitive' -> 'an intuitive')."
}
